$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Update data values
$ws.Range("B3").Value = 30
$ws.Range("D3").Value = 1145
$ws.Range("E3").Value = 921

$ws.Range("D4").Value = 7528
$ws.Range("E4").Value = 3274

# Recalculate formulas (F3/F4 are shared formulas referencing D/E)
$excel.Calculate()

# Update selection to E9
$ws.Range("E9").Select()

# Update window size/position (bookViews) to match the maximized-window
# state recorded in the saved workbook.
$win = $excel.ActiveWindow
$win.Left = -108
$win.Top = -108
$win.Width = 23256
$win.Height = 12720
